$d = $word.ActiveDocument

# Locate the paragraph containing "Ver no Jupiter Salvar em pdf Salvar em docx"
# by scanning the Paragraphs collection (robust against any pre-existing
# paragraph numbering quirks).
$count = $d.Paragraphs.Count
$verIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $verIdx = $i
        break
    }
}

if ($verIdx -eq -1) {
    Write-Output "Target paragraph not found; nothing changed."
} else {
    # Remove three whole paragraphs:
    #   verIdx - 1 : the blank paragraph right after the "LOB1039: ..." line
    #   verIdx     : "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   verIdx + 1 : the "(c) 2020 . Contact: ..." footer line
    $startPara = $d.Paragraphs.Item($verIdx - 1)
    $endPara = $d.Paragraphs.Item($verIdx + 1)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()

    Write-Output "Removed paragraphs $($verIdx - 1) through $($verIdx + 1)."
}
